# Natmi LR-pair results recomputed following Dr Hou's advice: the NATMI
# run now also includes the "ECs" sending cluster, so the table grows
# from 2 sending clusters x 4 target clusters (8 rows) to 3 x 4 (12 rows),
# and every numeric column is refreshed with the new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Full data table (rows 2-13), columns A-T.
# Column A: Sending cluster, B: Ligand symbol, C: Receptor symbol, D: Target cluster
$data = @(
    @("ECs","Tgfb2","Acvr1","ECs",2,0.6666666666666666,1.325336333333333,3.976009,0.02918077208126263,0.02918077208126263,3,1,5.436778333333334,16.310335,0.121853993972124,0.121853993972124,7.205559861446111,64.85003875301501,0.003555793625292102,0.003555793625292102),
    @("ECs","Tgfb2","Acvr1","FAPs",2,0.6666666666666666,1.325336333333333,3.976009,0.02918077208126263,0.02918077208126263,3,1,20.81735866666667,62.452076,0.4665774732677552,0.4665774732677551,27.59000180496489,248.310016244684,0.01361509090567777,0.01361509090567777),
    @("ECs","Tgfb2","Acvr1","M2",2,0.6666666666666666,1.325336333333333,3.976009,0.02918077208126263,0.02918077208126263,3,1,5.131792,15.395376,0.1150183643869107,0.1150183643869107,6.801350392709333,61.212153534384,0.003356324676334055,0.003356324676334055),
    @("ECs","Tgfb2","Acvr1","sCs",2,0.6666666666666666,1.325336333333333,3.976009,0.02918077208126263,0.02918077208126263,3,1,13.23122433333333,39.693673,0.2965501683732101,0.2965501683732101,17.53582234345077,157.822401091057,0.008653562873958701,0.008653562873958701),
    @("FAPs","Tgfb2","Acvr1","ECs",3,1,20.45485233333333,61.364557,0.4503674794711605,0.4503674794711605,3,1,5.436778333333334,16.310335,0.121853993972124,0.121853993972124,111.2084979773994,1000.876481796595,0.05487907612871948,0.05487907612871949),
    @("FAPs","Tgfb2","Acvr1","FAPs",3,1,20.45485233333333,61.364557,0.4503674794711605,0.4503674794711605,3,1,20.81735866666667,62.452076,0.4665774732677552,0.4665774732677551,425.8159974967035,3832.343977470332,0.2101313206136216,0.2101313206136216),
    @("FAPs","Tgfb2","Acvr1","M2",3,1,20.45485233333333,61.364557,0.4503674794711605,0.4503674794711605,3,1,5.131792,15.395376,0.1150183643869107,0.1150183643869107,104.9700475653813,944.7304280884319,0.05180053086182845,0.05180053086182845),
    @("FAPs","Tgfb2","Acvr1","sCs",3,1,20.45485233333333,61.364557,0.4503674794711605,0.4503674794711605,3,1,13.23122433333333,39.693673,0.2965501683732101,0.2965501683732101,270.64273992754,2435.784659347861,0.1335565518669909,0.1335565518669909),
    @("sCs","Tgfb2","Acvr1","ECs",3,1,23.63794933333334,70.913848,0.5204517484475769,0.5204517484475769,3,1,5.436778333333334,16.310335,0.121853993972124,0.121853993972124,128.5142907798978,1156.62861701908,0.06341912421811245,0.06341912421811245),
    @("sCs","Tgfb2","Acvr1","FAPs",3,1,23.63794933333334,70.913848,0.5204517484475769,0.5204517484475769,3,1,20.81735866666667,62.452076,0.4665774732677552,0.4665774732677551,492.0796694164943,4428.717024748448,0.2428310617484558,0.2428310617484557),
    @("sCs","Tgfb2","Acvr1","M2",3,1,23.63794933333334,70.913848,0.5204517484475769,0.5204517484475769,3,1,5.131792,15.395376,0.1150183643869107,0.1150183643869107,121.3050392852053,1091.745353566848,0.05986150884874818,0.05986150884874818),
    @("sCs","Tgfb2","Acvr1","sCs",3,1,23.63794933333334,70.913848,0.5204517484475769,0.5204517484475769,3,1,13.23122433333333,39.693673,0.2965501683732101,0.2965501683732101,312.7590104093004,2814.831093683704,0.1543400536322605,0.1543400536322605)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value = $rowVals[$c]
    }
}
